$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Round the coordinate values in row 5 to whole numbers
$ws.Range("Q5").Value = 450824
$ws.Range("R5").Value = 7031623

# Clear the "Starttid" (Z5) and "Sluttid" (AB5) time cells for row 5
$ws.Range("Z5").ClearContents()
$ws.Range("AB5").ClearContents()
